$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(7, 8).Value = 0
$ws.Cells.Item(7, 10).Value = 0
$ws.Cells.Item(7, 12).Value = 0
$ws.Cells.Item(7, 14).Value = $null

$ws.Cells.Item(10, 8).Value = 23050
$ws.Cells.Item(10, 10).Value = 23050
$ws.Cells.Item(10, 12).Value = 23050
$ws.Cells.Item(10, 14).Value = -23636

$ws.Cells.Item(14, 8).Value = 0
$ws.Cells.Item(14, 10).Value = 0
$ws.Cells.Item(14, 12).Value = 0
$ws.Cells.Item(14, 14).Value = $null

$ws.Cells.Item(39, 8).Value = 238.72728
$ws.Cells.Item(39, 9).Value = 91.77778000000001
$ws.Cells.Item(39, 11).Value = 275.33334
$ws.Cells.Item(39, 13).Value = 20.66665999999998

$ws.Cells.Item(40, 8).Value = 1211.375
$ws.Cells.Item(40, 10).Value = 1211.375
$ws.Cells.Item(40, 12).Value = 1211.375
$ws.Cells.Item(40, 14).Value = -1561.375

$ws.Cells.Item(43, 8).Value = 21199.2
$ws.Cells.Item(43, 10).Value = 1599
$ws.Cells.Item(43, 12).Value = 1599
$ws.Cells.Item(43, 14).Value = -1737

$ws.Cells.Item(54, 8).Value = 8703.666999999999
$ws.Cells.Item(54, 9).Value = 8703.666999999999
$ws.Cells.Item(54, 11).Value = 8703.666999999999
$ws.Cells.Item(54, 13).Value = -8217.666999999999

$ws.Cells.Item(62, 8).Value = 6676400
$ws.Cells.Item(62, 9).Value = 13340800
$ws.Cells.Item(62, 11).Value = 13340800
$ws.Cells.Item(62, 13).Value = -13340176

$ws.Cells.Item(64, 8).Value = 9601.25
$ws.Cells.Item(64, 10).Value = 11621.444
$ws.Cells.Item(64, 12).Value = 11621.444
$ws.Cells.Item(64, 14).Value = -12117.444

$ws.Cells.Item(65, 8).Value = 6676400
$ws.Cells.Item(65, 9).Value = 13340800
$ws.Cells.Item(65, 11).Value = 66704000
$ws.Cells.Item(65, 13).Value = -66700880

$ws.Cells.Item(67, 8).Value = 9601.25
$ws.Cells.Item(67, 10).Value = 11621.444
$ws.Cells.Item(67, 12).Value = 11621.444
$ws.Cells.Item(67, 14).Value = -13337.444

$ws.Cells.Item(76, 8).Value = 5514.8
$ws.Cells.Item(76, 9).Value = 5191.6665
$ws.Cells.Item(76, 11).Value = 5191.6665
$ws.Cells.Item(76, 13).Value = -4876.6665

$ws.Cells.Item(79, 8).Value = 5514.8
$ws.Cells.Item(79, 9).Value = 5191.6665
$ws.Cells.Item(79, 11).Value = 5191.6665
$ws.Cells.Item(79, 13).Value = -4099.6665

$ws.Cells.Item(86, 8).Value = 10000
$ws.Cells.Item(86, 9).Value = 0
$ws.Cells.Item(86, 10).Value = 10000
$ws.Cells.Item(86, 11).Value = 0
$ws.Cells.Item(86, 12).Value = 10000
$ws.Cells.Item(86, 14).Value = -12246
$ws.Cells.Item(86, 13).Value = $null

$ws.Cells.Item(89, 8).Value = 10000
$ws.Cells.Item(89, 9).Value = 0
$ws.Cells.Item(89, 10).Value = 10000
$ws.Cells.Item(89, 11).Value = 0
$ws.Cells.Item(89, 12).Value = 50000
$ws.Cells.Item(89, 14).Value = -61232
$ws.Cells.Item(89, 13).Value = $null

$ws.Cells.Item(96, 8).Value = 890.7826
$ws.Cells.Item(96, 9).Value = 1018.125
$ws.Cells.Item(96, 10).Value = 599.7143
$ws.Cells.Item(96, 11).Value = 3054.375
$ws.Cells.Item(96, 12).Value = 1799.1429
$ws.Cells.Item(96, 13).Value = -1681.375
$ws.Cells.Item(96, 14).Value = -4545.1429

$ws.Cells.Item(106, 8).Value = 50001660
$ws.Cells.Item(106, 10).Value = 0
$ws.Cells.Item(106, 12).Value = 0
$ws.Cells.Item(106, 14).Value = $null

$ws.Cells.Item(113, 8).Value = 55363.41
$ws.Cells.Item(113, 10).Value = 11780.818
$ws.Cells.Item(113, 12).Value = 11780.818
$ws.Cells.Item(113, 14).Value = -18288.818

$ws.Cells.Item(129, 8).Value = 3294.75
$ws.Cells.Item(129, 9).Value = 3294.75
$ws.Cells.Item(129, 11).Value = 9884.25
$ws.Cells.Item(129, 13).Value = -4884.25

$ws.Cells.Item(130, 8).Value = 150000
$ws.Cells.Item(130, 10).Value = 150000
$ws.Cells.Item(130, 12).Value = 150000
$ws.Cells.Item(130, 14).Value = -160040

$ws.Cells.Item(131, 8).Value = 9909.277
$ws.Cells.Item(131, 9).Value = 2250.7058
$ws.Cells.Item(131, 11).Value = 6752.117400000001
$ws.Cells.Item(131, 13).Value = -1712.117400000001

$ws.Cells.Item(132, 8).Value = 9805649
$ws.Cells.Item(132, 9).Value = 11906212
$ws.Cells.Item(132, 11).Value = 35718636
$ws.Cells.Item(132, 13).Value = -35716106

$ws.Cells.Item(137, 8).Value = 1704.6666
$ws.Cells.Item(137, 9).Value = 1612.625
$ws.Cells.Item(137, 10).Value = 1999.2
$ws.Cells.Item(137, 11).Value = 4837.875
$ws.Cells.Item(137, 12).Value = 5997.6
$ws.Cells.Item(137, 13).Value = -2287.875
$ws.Cells.Item(137, 14).Value = -11097.6

$ws.Cells.Item(138, 8).Value = 3322.3247
$ws.Cells.Item(138, 9).Value = 1814.2667
$ws.Cells.Item(138, 10).Value = 4284.915
$ws.Cells.Item(138, 11).Value = 5442.800099999999
$ws.Cells.Item(138, 12).Value = 12854.745
$ws.Cells.Item(138, 13).Value = -302.8000999999995
$ws.Cells.Item(138, 14).Value = -23134.745

$ws.Cells.Item(141, 8).Value = 1863.3334
$ws.Cells.Item(141, 9).Value = 1706.5
$ws.Cells.Item(141, 11).Value = 5119.5
$ws.Cells.Item(141, 13).Value = 60.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(22, 8).Value = 8476.875
$ws.Cells.Item(22, 9).Value = 4272
$ws.Cells.Item(22, 10).Value = 10999.8
$ws.Cells.Item(22, 11).Value = 4272
$ws.Cells.Item(22, 12).Value = 10999.8
$ws.Cells.Item(22, 13).Value = -3973
$ws.Cells.Item(22, 14).Value = -11597.8

$ws.Cells.Item(32, 8).Value = 3728.147
$ws.Cells.Item(32, 9).Value = 2848.8196
$ws.Cells.Item(32, 11).Value = 2848.8196
$ws.Cells.Item(32, 13).Value = -2561.8196

$ws.Cells.Item(46, 8).Value = 3166.3333
$ws.Cells.Item(46, 10).Value = 3166.3333
$ws.Cells.Item(46, 12).Value = 3166.3333
$ws.Cells.Item(46, 14).Value = -3804.3333

$ws.Cells.Item(74, 8).Value = 5046.378
$ws.Cells.Item(74, 9).Value = 1029.7949
$ws.Cells.Item(74, 10).Value = 31154.166
$ws.Cells.Item(74, 11).Value = 1029.7949
$ws.Cells.Item(74, 12).Value = 31154.166
$ws.Cells.Item(74, 13).Value = -155.7949000000001
$ws.Cells.Item(74, 14).Value = -32902.166

$ws.Cells.Item(77, 8).Value = 5046.378
$ws.Cells.Item(77, 9).Value = 1029.7949
$ws.Cells.Item(77, 10).Value = 31154.166
$ws.Cells.Item(77, 11).Value = 5148.9745
$ws.Cells.Item(77, 12).Value = 155770.83
$ws.Cells.Item(77, 13).Value = -780.9745000000003
$ws.Cells.Item(77, 14).Value = -164506.83

$ws.Cells.Item(86, 8).Value = 36051.668
$ws.Cells.Item(86, 9).Value = 18000
$ws.Cells.Item(86, 10).Value = 39662
$ws.Cells.Item(86, 11).Value = 18000
$ws.Cells.Item(86, 12).Value = 39662
$ws.Cells.Item(86, 14).Value = -42034
$ws.Cells.Item(86, 13).Value = -16814

$ws.Cells.Item(89, 8).Value = 36051.668
$ws.Cells.Item(89, 9).Value = 18000
$ws.Cells.Item(89, 10).Value = 39662
$ws.Cells.Item(89, 11).Value = 54000
$ws.Cells.Item(89, 12).Value = 118986
$ws.Cells.Item(89, 14).Value = -130842
$ws.Cells.Item(89, 13).Value = -48072

$ws.Cells.Item(97, 8).Value = 345.18182
$ws.Cells.Item(97, 9).Value = 395.66666
$ws.Cells.Item(97, 10).Value = 118
$ws.Cells.Item(97, 11).Value = 395.66666
$ws.Cells.Item(97, 12).Value = 118
$ws.Cells.Item(97, 13).Value = 100.33334
$ws.Cells.Item(97, 14).Value = -1110

$ws.Cells.Item(102, 8).Value = 4928.294
$ws.Cells.Item(102, 10).Value = 6228.6665
$ws.Cells.Item(102, 12).Value = 6228.6665
$ws.Cells.Item(102, 14).Value = -9472.666499999999

$ws.Cells.Item(122, 8).Value = 2088
$ws.Cells.Item(122, 9).Value = 1865.2727
$ws.Cells.Item(122, 10).Value = 2904.6667
$ws.Cells.Item(122, 11).Value = 5595.8181
$ws.Cells.Item(122, 12).Value = 8714.000100000001
$ws.Cells.Item(122, 13).Value = -3145.8181
$ws.Cells.Item(122, 14).Value = -13614.0001

$ws.Cells.Item(131, 8).Value = 116666.664
$ws.Cells.Item(131, 10).Value = 116666.664
$ws.Cells.Item(131, 12).Value = 116666.664
$ws.Cells.Item(131, 14).Value = -126746.664

$ws.Cells.Item(132, 8).Value = 5630.5
$ws.Cells.Item(132, 9).Value = 5559.8
$ws.Cells.Item(132, 11).Value = 16679.4
$ws.Cells.Item(132, 13).Value = -14149.4

$ws.Cells.Item(141, 8).Value = 36851.668
$ws.Cells.Item(141, 9).Value = 15000
$ws.Cells.Item(141, 11).Value = 15000
$ws.Cells.Item(141, 13).Value = -9820

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(25, 8).Value = 1810.3334
$ws.Cells.Item(25, 10).Value = 2338.6667
$ws.Cells.Item(25, 12).Value = 2338.6667
$ws.Cells.Item(25, 14).Value = -2808.6667

$ws.Cells.Item(38, 8).Value = 49000
$ws.Cells.Item(38, 10).Value = 49000
$ws.Cells.Item(38, 12).Value = 49000
$ws.Cells.Item(38, 14).Value = -49832

$ws.Cells.Item(86, 8).Value = 4421.0303
$ws.Cells.Item(86, 9).Value = 1522.8518
$ws.Cells.Item(86, 10).Value = 17462.834
$ws.Cells.Item(86, 11).Value = 1522.8518
$ws.Cells.Item(86, 12).Value = 17462.834
$ws.Cells.Item(86, 13).Value = -399.8517999999999
$ws.Cells.Item(86, 14).Value = -19708.834

$ws.Cells.Item(89, 8).Value = 4421.0303
$ws.Cells.Item(89, 9).Value = 1522.8518
$ws.Cells.Item(89, 10).Value = 17462.834
$ws.Cells.Item(89, 11).Value = 7614.259
$ws.Cells.Item(89, 12).Value = 87314.17
$ws.Cells.Item(89, 13).Value = -1998.259
$ws.Cells.Item(89, 14).Value = -98546.17

$ws.Cells.Item(94, 8).Value = 3008.4285
$ws.Cells.Item(94, 9).Value = 3250
$ws.Cells.Item(94, 11).Value = 3250
$ws.Cells.Item(94, 13).Value = -2799

$ws.Cells.Item(99, 8).Value = 4333
$ws.Cells.Item(99, 9).Value = 5000
$ws.Cells.Item(99, 10).Value = 3999.5
$ws.Cells.Item(99, 11).Value = 5000
$ws.Cells.Item(99, 12).Value = 3999.5
$ws.Cells.Item(99, 13).Value = -3502
$ws.Cells.Item(99, 14).Value = -6995.5

$ws.Cells.Item(105, 8).Value = 1361.2778
$ws.Cells.Item(105, 9).Value = 1218.9375
$ws.Cells.Item(105, 11).Value = 1218.9375
$ws.Cells.Item(105, 13).Value = 528.0625

$ws.Cells.Item(134, 8).Value = 1810.65
$ws.Cells.Item(134, 9).Value = 1847.7297
$ws.Cells.Item(134, 10).Value = 1353.3334
$ws.Cells.Item(134, 11).Value = 5543.189100000001
$ws.Cells.Item(134, 12).Value = 4060.0002
$ws.Cells.Item(134, 13).Value = -3008.189100000001
$ws.Cells.Item(134, 14).Value = -9130.0002

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(3, 8).Value = 5455.7144
$ws.Cells.Item(3, 9).Value = 3480.6667
$ws.Cells.Item(3, 10).Value = 6937
$ws.Cells.Item(3, 11).Value = 3480.6667
$ws.Cells.Item(3, 12).Value = 6937
$ws.Cells.Item(3, 13).Value = -3367.6667
$ws.Cells.Item(3, 14).Value = -7163

$ws.Cells.Item(20, 8).Value = 250000
$ws.Cells.Item(20, 10).Value = 250000
$ws.Cells.Item(20, 12).Value = 250000
$ws.Cells.Item(20, 14).Value = -250472

$ws.Cells.Item(30, 8).Value = 250000
$ws.Cells.Item(30, 10).Value = 250000
$ws.Cells.Item(30, 12).Value = 250000
$ws.Cells.Item(30, 14).Value = -250182

$ws.Cells.Item(31, 8).Value = 67696.31
$ws.Cells.Item(31, 9).Value = 78744
$ws.Cells.Item(31, 11).Value = 78744
$ws.Cells.Item(31, 13).Value = -78449

$ws.Cells.Item(33, 8).Value = 24000.125
$ws.Cells.Item(33, 9).Value = 504
$ws.Cells.Item(33, 11).Value = 504
$ws.Cells.Item(33, 13).Value = -125

$ws.Cells.Item(34, 8).Value = 67696.31
$ws.Cells.Item(34, 9).Value = 78744
$ws.Cells.Item(34, 11).Value = 78744
$ws.Cells.Item(34, 13).Value = -78542

$ws.Cells.Item(105, 8).Value = 1166.5
$ws.Cells.Item(105, 10).Value = 0
$ws.Cells.Item(105, 12).Value = 0
$ws.Cells.Item(105, 14).Value = $null

$ws.Cells.Item(128, 8).Value = 250000
$ws.Cells.Item(128, 10).Value = 250000
$ws.Cells.Item(128, 12).Value = 250000
$ws.Cells.Item(128, 14).Value = -259960

$ws.Cells.Item(132, 8).Value = 4385.926
$ws.Cells.Item(132, 9).Value = 4312.8
$ws.Cells.Item(132, 11).Value = 12938.4
$ws.Cells.Item(132, 13).Value = -10408.4

$ws.Cells.Item(134, 8).Value = 16991.572
$ws.Cells.Item(134, 9).Value = 7671.273
$ws.Cells.Item(134, 11).Value = 23013.819
$ws.Cells.Item(134, 13).Value = -20478.819

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(3, 8).Value = 12833.333
$ws.Cells.Item(3, 9).Value = 11000
$ws.Cells.Item(3, 10).Value = 13750
$ws.Cells.Item(3, 11).Value = 33000
$ws.Cells.Item(3, 12).Value = 41250
$ws.Cells.Item(3, 13).Value = -32888
$ws.Cells.Item(3, 14).Value = -41474

$ws.Cells.Item(5, 8).Value = 11359.6
$ws.Cells.Item(5, 9).Value = 0
$ws.Cells.Item(5, 11).Value = 0
$ws.Cells.Item(5, 13).Value = $null

$ws.Cells.Item(7, 8).Value = 544.3333
$ws.Cells.Item(7, 10).Value = 100
$ws.Cells.Item(7, 12).Value = 300
$ws.Cells.Item(7, 14).Value = -524

$ws.Cells.Item(19, 8).Value = 799
$ws.Cells.Item(19, 9).Value = 799
$ws.Cells.Item(19, 10).Value = 799
$ws.Cells.Item(19, 11).Value = 2397
$ws.Cells.Item(19, 12).Value = 2397
$ws.Cells.Item(19, 13).Value = -2223
$ws.Cells.Item(19, 14).Value = -2745

$ws.Cells.Item(22, 8).Value = 4653.385
$ws.Cells.Item(22, 9).Value = 400
$ws.Cells.Item(22, 11).Value = 1200
$ws.Cells.Item(22, 13).Value = -1031

$ws.Cells.Item(27, 8).Value = 4653.385
$ws.Cells.Item(27, 9).Value = 400
$ws.Cells.Item(27, 11).Value = 1200
$ws.Cells.Item(27, 13).Value = -1098

$ws.Cells.Item(40, 8).Value = 323.9
$ws.Cells.Item(40, 9).Value = 189.85715
$ws.Cells.Item(40, 10).Value = 636.6667
$ws.Cells.Item(40, 11).Value = 759.4286
$ws.Cells.Item(40, 12).Value = 2546.6668
$ws.Cells.Item(40, 13).Value = -690.4286
$ws.Cells.Item(40, 14).Value = -2684.6668

$ws.Cells.Item(47, 8).Value = 88.25
$ws.Cells.Item(47, 9).Value = 88.25
$ws.Cells.Item(47, 11).Value = 264.75
$ws.Cells.Item(47, 13).Value = 166.25

$ws.Cells.Item(80, 8).Value = 8332.333000000001
$ws.Cells.Item(80, 9).Value = 4999
$ws.Cells.Item(80, 10).Value = 9999
$ws.Cells.Item(80, 11).Value = 14997
$ws.Cells.Item(80, 12).Value = 29997
$ws.Cells.Item(80, 13).Value = -14061
$ws.Cells.Item(80, 14).Value = -31869

$ws.Cells.Item(81, 8).Value = 4425.8184
$ws.Cells.Item(81, 9).Value = 986
$ws.Cells.Item(81, 10).Value = 5190.222
$ws.Cells.Item(81, 11).Value = 2958
$ws.Cells.Item(81, 12).Value = 15570.666
$ws.Cells.Item(81, 13).Value = -1835
$ws.Cells.Item(81, 14).Value = -17816.666

$ws.Cells.Item(83, 8).Value = 8332.333000000001
$ws.Cells.Item(83, 9).Value = 4999
$ws.Cells.Item(83, 10).Value = 9999
$ws.Cells.Item(83, 11).Value = 44991
$ws.Cells.Item(83, 12).Value = 89991
$ws.Cells.Item(83, 13).Value = -40311
$ws.Cells.Item(83, 14).Value = -99351

$ws.Cells.Item(84, 8).Value = 4425.8184
$ws.Cells.Item(84, 9).Value = 986
$ws.Cells.Item(84, 10).Value = 5190.222
$ws.Cells.Item(84, 11).Value = 8874
$ws.Cells.Item(84, 12).Value = 46711.998
$ws.Cells.Item(84, 13).Value = -3258
$ws.Cells.Item(84, 14).Value = -57943.998

$ws.Cells.Item(109, 8).Value = 462.5
$ws.Cells.Item(109, 9).Value = 462.5
$ws.Cells.Item(109, 11).Value = 1387.5
$ws.Cells.Item(109, 13).Value = -347.5

$ws.Cells.Item(135, 8).Value = 11359.6
$ws.Cells.Item(135, 9).Value = 0
$ws.Cells.Item(135, 11).Value = 0
$ws.Cells.Item(135, 13).Value = $null

$ws.Cells.Item(140, 8).Value = 2844.2727
$ws.Cells.Item(140, 9).Value = 2844.2727
$ws.Cells.Item(140, 10).Value = 0
$ws.Cells.Item(140, 11).Value = 8532.8181
$ws.Cells.Item(140, 12).Value = 0
$ws.Cells.Item(140, 13).Value = -3352.8181
$ws.Cells.Item(140, 14).Value = $null

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(5, 8).Value = 7168.3335
$ws.Cells.Item(5, 9).Value = 752.5
$ws.Cells.Item(5, 10).Value = 20000
$ws.Cells.Item(5, 11).Value = 752.5
$ws.Cells.Item(5, 12).Value = 20000
$ws.Cells.Item(5, 13).Value = -640.5
$ws.Cells.Item(5, 14).Value = -20224

$ws.Cells.Item(15, 8).Value = 40000
$ws.Cells.Item(15, 10).Value = 40000
$ws.Cells.Item(15, 12).Value = 40000
$ws.Cells.Item(15, 14).Value = -40576

$ws.Cells.Item(41, 8).Value = 3000
$ws.Cells.Item(41, 9).Value = 3000
$ws.Cells.Item(41, 11).Value = 3000
$ws.Cells.Item(41, 13).Value = -2645

$ws.Cells.Item(48, 8).Value = 0
$ws.Cells.Item(48, 10).Value = 0
$ws.Cells.Item(48, 12).Value = 0
$ws.Cells.Item(48, 14).Value = $null

$ws.Cells.Item(80, 8).Value = 8752.25
$ws.Cells.Item(80, 9).Value = 4005
$ws.Cells.Item(80, 10).Value = 10334.667
$ws.Cells.Item(80, 11).Value = 4005
$ws.Cells.Item(80, 12).Value = 10334.667
$ws.Cells.Item(80, 14).Value = -12330.667
$ws.Cells.Item(80, 13).Value = -3007

$ws.Cells.Item(81, 8).Value = 40000
$ws.Cells.Item(81, 10).Value = 40000
$ws.Cells.Item(81, 12).Value = 40000
$ws.Cells.Item(81, 14).Value = -41996

$ws.Cells.Item(83, 8).Value = 8752.25
$ws.Cells.Item(83, 9).Value = 4005
$ws.Cells.Item(83, 10).Value = 10334.667
$ws.Cells.Item(83, 11).Value = 20025
$ws.Cells.Item(83, 12).Value = 51673.335
$ws.Cells.Item(83, 14).Value = -61657.335
$ws.Cells.Item(83, 13).Value = -15033

$ws.Cells.Item(84, 8).Value = 40000
$ws.Cells.Item(84, 10).Value = 40000
$ws.Cells.Item(84, 12).Value = 120000
$ws.Cells.Item(84, 14).Value = -129984

$ws.Cells.Item(97, 8).Value = 1544.7931
$ws.Cells.Item(97, 9).Value = 960.4545000000001
$ws.Cells.Item(97, 10).Value = 3381.2856
$ws.Cells.Item(97, 11).Value = 960.4545000000001
$ws.Cells.Item(97, 12).Value = 3381.2856
$ws.Cells.Item(97, 13).Value = -464.4545000000001
$ws.Cells.Item(97, 14).Value = -4373.2856

$ws.Cells.Item(104, 8).Value = 31000
$ws.Cells.Item(104, 10).Value = 31000
$ws.Cells.Item(104, 12).Value = 31000
$ws.Cells.Item(104, 14).Value = -37988

$ws.Cells.Item(113, 8).Value = 2116.3333
$ws.Cells.Item(113, 10).Value = 2650
$ws.Cells.Item(113, 12).Value = 2650
$ws.Cells.Item(113, 14).Value = -6990

$ws.Cells.Item(126, 8).Value = 21202.445
$ws.Cells.Item(126, 9).Value = 38370.777
$ws.Cells.Item(126, 10).Value = 4034.111
$ws.Cells.Item(126, 11).Value = 115112.331
$ws.Cells.Item(126, 12).Value = 12102.333
$ws.Cells.Item(126, 13).Value = -112642.331
$ws.Cells.Item(126, 14).Value = -17042.333

$ws.Cells.Item(132, 8).Value = 5579.9
$ws.Cells.Item(132, 9).Value = 4925
$ws.Cells.Item(132, 10).Value = 6016.5
$ws.Cells.Item(132, 11).Value = 14775
$ws.Cells.Item(132, 12).Value = 18049.5
$ws.Cells.Item(132, 13).Value = -12245
$ws.Cells.Item(132, 14).Value = -23109.5

$ws.Cells.Item(134, 8).Value = 74826
$ws.Cells.Item(134, 10).Value = 74826
$ws.Cells.Item(134, 12).Value = 224478
$ws.Cells.Item(134, 14).Value = -229548

$ws.Cells.Item(136, 8).Value = 35777.723
$ws.Cells.Item(136, 10).Value = 35777.723
$ws.Cells.Item(136, 12).Value = 107333.169
$ws.Cells.Item(136, 14).Value = -112433.169

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(22, 8).Value = 1505
$ws.Cells.Item(22, 9).Value = 1800
$ws.Cells.Item(22, 10).Value = 1446
$ws.Cells.Item(22, 11).Value = 1800
$ws.Cells.Item(22, 12).Value = 1446
$ws.Cells.Item(22, 13).Value = -1505
$ws.Cells.Item(22, 14).Value = -2036

$ws.Cells.Item(27, 8).Value = 1505
$ws.Cells.Item(27, 9).Value = 1800
$ws.Cells.Item(27, 10).Value = 1446
$ws.Cells.Item(27, 11).Value = 1800
$ws.Cells.Item(27, 12).Value = 1446
$ws.Cells.Item(27, 13).Value = -1693
$ws.Cells.Item(27, 14).Value = -1660

$ws.Cells.Item(40, 8).Value = 5277.9287
$ws.Cells.Item(40, 9).Value = 3543.6667
$ws.Cells.Item(40, 11).Value = 3543.6667
$ws.Cells.Item(40, 13).Value = -3407.6667

$ws.Cells.Item(42, 8).Value = 12853.214

$ws.Cells.Item(43, 8).Value = 2520000
$ws.Cells.Item(43, 9).Value = 7500000
$ws.Cells.Item(43, 11).Value = 7500000
$ws.Cells.Item(43, 13).Value = -7499807

$ws.Cells.Item(46, 8).Value = 2704.6924
$ws.Cells.Item(46, 9).Value = 1341.2858
$ws.Cells.Item(46, 10).Value = 4295.3335
$ws.Cells.Item(46, 11).Value = 1341.2858
$ws.Cells.Item(46, 12).Value = 4295.3335
$ws.Cells.Item(46, 13).Value = -1153.2858
$ws.Cells.Item(46, 14).Value = -4671.3335

$ws.Cells.Item(49, 8).Value = 12853.214

$ws.Cells.Item(55, 8).Value = 942.3333
$ws.Cells.Item(55, 9).Value = 925
$ws.Cells.Item(55, 10).Value = 951
$ws.Cells.Item(55, 11).Value = 925
$ws.Cells.Item(55, 12).Value = 951
$ws.Cells.Item(55, 13).Value = -752
$ws.Cells.Item(55, 14).Value = -1297

$ws.Cells.Item(61, 8).Value = 3977.0625
$ws.Cells.Item(61, 9).Value = 3977.0625
$ws.Cells.Item(61, 11).Value = 3977.0625
$ws.Cells.Item(61, 13).Value = -3775.0625

$ws.Cells.Item(68, 8).Value = 4111
$ws.Cells.Item(68, 10).Value = 5222
$ws.Cells.Item(68, 12).Value = 5222
$ws.Cells.Item(68, 14).Value = -6720

$ws.Cells.Item(71, 8).Value = 4111
$ws.Cells.Item(71, 10).Value = 5222
$ws.Cells.Item(71, 12).Value = 26110
$ws.Cells.Item(71, 14).Value = -33598

$ws.Cells.Item(113, 8).Value = 3977.0625
$ws.Cells.Item(113, 9).Value = 3977.0625
$ws.Cells.Item(113, 11).Value = 3977.0625
$ws.Cells.Item(113, 13).Value = -1807.0625

$ws.Cells.Item(115, 8).Value = 0
$ws.Cells.Item(115, 10).Value = 0
$ws.Cells.Item(115, 12).Value = 0
$ws.Cells.Item(115, 14).Value = $null

$ws.Cells.Item(122, 8).Value = 4024.5476
$ws.Cells.Item(122, 9).Value = 3481.1924
$ws.Cells.Item(122, 10).Value = 4907.5
$ws.Cells.Item(122, 11).Value = 10443.5772
$ws.Cells.Item(122, 12).Value = 14722.5
$ws.Cells.Item(122, 13).Value = -7993.5772
$ws.Cells.Item(122, 14).Value = -19622.5

$ws.Cells.Item(131, 8).Value = 28200
$ws.Cells.Item(131, 9).Value = 34000
$ws.Cells.Item(131, 11).Value = 34000
$ws.Cells.Item(131, 13).Value = -28960

$ws.Cells.Item(132, 8).Value = 5671.4287
$ws.Cells.Item(132, 9).Value = 4925
$ws.Cells.Item(132, 10).Value = 6666.6665
$ws.Cells.Item(132, 11).Value = 14775
$ws.Cells.Item(132, 12).Value = 19999.9995
$ws.Cells.Item(132, 13).Value = -12245
$ws.Cells.Item(132, 14).Value = -25059.9995

$ws.Cells.Item(136, 8).Value = 4971.107
$ws.Cells.Item(136, 9).Value = 4312.4346
$ws.Cells.Item(136, 10).Value = 8001
$ws.Cells.Item(136, 11).Value = 12937.3038
$ws.Cells.Item(136, 12).Value = 24003
$ws.Cells.Item(136, 13).Value = -10387.3038
$ws.Cells.Item(136, 14).Value = -29103

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(18, 8).Value = 2492.5881
$ws.Cells.Item(18, 9).Value = 899.75
$ws.Cells.Item(18, 10).Value = 2982.6924
$ws.Cells.Item(18, 11).Value = 899.75
$ws.Cells.Item(18, 12).Value = 2982.6924
$ws.Cells.Item(18, 13).Value = -726.75
$ws.Cells.Item(18, 14).Value = -3328.6924

$ws.Cells.Item(29, 8).Value = 10500
$ws.Cells.Item(29, 9).Value = 11250
$ws.Cells.Item(29, 10).Value = 9000
$ws.Cells.Item(29, 11).Value = 11250
$ws.Cells.Item(29, 12).Value = 9000
$ws.Cells.Item(29, 13).Value = -10960
$ws.Cells.Item(29, 14).Value = -9580

$ws.Cells.Item(39, 8).Value = 9994
$ws.Cells.Item(39, 9).Value = 9994
$ws.Cells.Item(39, 11).Value = 9994
$ws.Cells.Item(39, 13).Value = -9581

$ws.Cells.Item(81, 8).Value = 8140.9546
$ws.Cells.Item(81, 9).Value = 20692.2
$ws.Cells.Item(81, 11).Value = 41384.4
$ws.Cells.Item(81, 13).Value = -40323.4

$ws.Cells.Item(84, 8).Value = 8140.9546
$ws.Cells.Item(84, 9).Value = 20692.2
$ws.Cells.Item(84, 11).Value = 206922
$ws.Cells.Item(84, 13).Value = -201618

$ws.Cells.Item(100, 8).Value = 1496.8182
$ws.Cells.Item(100, 9).Value = 1103.3334
$ws.Cells.Item(100, 11).Value = 2206.6668
$ws.Cells.Item(100, 13).Value = -1665.6668

$ws.Cells.Item(123, 8).Value = 41000
$ws.Cells.Item(123, 10).Value = 41000
$ws.Cells.Item(123, 12).Value = 41000
$ws.Cells.Item(123, 14).Value = -50800

$ws.Cells.Item(132, 8).Value = 7684.231
$ws.Cells.Item(132, 9).Value = 7487.375
$ws.Cells.Item(132, 10).Value = 7999.2
$ws.Cells.Item(132, 11).Value = 22462.125
$ws.Cells.Item(132, 12).Value = 23997.6
$ws.Cells.Item(132, 13).Value = -19932.125
$ws.Cells.Item(132, 14).Value = -29057.6

$ws.Cells.Item(136, 8).Value = 3107.5625
$ws.Cells.Item(136, 9).Value = 2981.4
$ws.Cells.Item(136, 10).Value = 5000
$ws.Cells.Item(136, 11).Value = 8944.200000000001
$ws.Cells.Item(136, 12).Value = 15000
$ws.Cells.Item(136, 13).Value = -6394.200000000001
$ws.Cells.Item(136, 14).Value = -20100
